$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2042857142857143
$ws.Range("C2").Value = 0.5457142857142857
$ws.Range("J2").Value = 0.008571428571428572
$ws.Range("O2").Value = 0.001428571428571429
$ws.Range("P2").Value = 0.1371428571428571
$ws.Range("S2").Value = 0.1028571428571429
$ws.Range("B3").Value = 0.0100250626566416
$ws.Range("C3").Value = 0.02005012531328321
$ws.Range("J3").Value = 0.02255639097744361
$ws.Range("P3").Value = 0.7293233082706767
$ws.Range("S3").Value = 0.2180451127819549
$ws.Range("J4").Value = 0.04123711340206185
$ws.Range("O4").Value = 0.01030927835051546
$ws.Range("P4").Value = 0.7319587628865979
$ws.Range("S4").Value = 0.2164948453608248
$ws.Range("B6").Value = 0.08918406072106262
$ws.Range("D6").Value = 0.01328273244781784
$ws.Range("F6").Value = 0.08728652751423149
$ws.Range("J6").Value = 0.1916508538899431
$ws.Range("O6").Value = 0.04174573055028463
$ws.Range("Q6").Value = 0.1555977229601518
$ws.Range("R6").Value = 0.08159392789373814
$ws.Range("S6").Value = 0.3396584440227704
$ws.Range("B7").Value = 0.1328502415458937
$ws.Range("D7").Value = 0.02898550724637681
$ws.Range("F7").Value = 0.07004830917874397
$ws.Range("J7").Value = 0.09420289855072464
$ws.Range("O7").Value = 0.01449275362318841
$ws.Range("Q7").Value = 0.1811594202898551
$ws.Range("R7").Value = 0.0748792270531401
$ws.Range("S7").Value = 0.4033816425120773
$ws.Range("B8").Value = 0.1
$ws.Range("D8").Value = 0.01568627450980392
$ws.Range("F8").Value = 0.07156862745098039
$ws.Range("J8").Value = 0.1176470588235294
$ws.Range("O8").Value = 0.02156862745098039
$ws.Range("Q8").Value = 0.1450980392156863
$ws.Range("R8").Value = 0.08235294117647059
$ws.Range("S8").Value = 0.446078431372549
$ws.Range("B9").Value = 0.1214689265536723
$ws.Range("D9").Value = 0.01412429378531073
$ws.Range("E9").Value = 0.002824858757062147
$ws.Range("F9").Value = 0.08757062146892655
$ws.Range("J9").Value = 0.1129943502824859
$ws.Range("O9").Value = 0.02824858757062147
$ws.Range("Q9").Value = 0.1327683615819209
$ws.Range("R9").Value = 0.08757062146892655
$ws.Range("S9").Value = 0.4124293785310734
$ws.Range("B10").Value = 0.1179952644041042
$ws.Range("D10").Value = 0.02407261247040253
$ws.Range("E10").Value = 0.0007892659826361484
$ws.Range("F10").Value = 0.07537490134175218
$ws.Range("J10").Value = 0.1152328334648777
$ws.Range("O10").Value = 0.01894238358326756
$ws.Range("Q10").Value = 0.1973164956590371
$ws.Range("R10").Value = 0.07576953433307025
$ws.Range("S10").Value = 0.3745067087608524
$ws.Range("G11").Value = 0.1409495548961424
$ws.Range("J11").Value = 0.1008902077151335
$ws.Range("K11").Value = 0.2270029673590505
$ws.Range("L11").Value = 0.5178041543026706
$ws.Range("S11").Value = 0.01335311572700297
$ws.Range("G12").Value = 0.7211796246648794
$ws.Range("J12").Value = 0.1957104557640751
$ws.Range("K12").Value = 0.01072386058981233
$ws.Range("L12").Value = 0.04021447721179625
$ws.Range("S12").Value = 0.032171581769437
$ws.Range("G13").Value = 0.6931818181818182
$ws.Range("J13").Value = 0.25
$ws.Range("S13").Value = 0.05681818181818182
$ws.Range("G14").Value = 0.8571428571428571
$ws.Range("J14").Value = 0.1428571428571428
$ws.Range("F15").Value = 0.02444444444444445
$ws.Range("H15").Value = 0.1533333333333333
$ws.Range("I15").Value = 0.06444444444444444
$ws.Range("J15").Value = 0.3444444444444444
$ws.Range("K15").Value = 0.06444444444444444
$ws.Range("M15").Value = 0.01333333333333333
$ws.Range("N15").Value = 0.002222222222222222
$ws.Range("O15").Value = 0.05777777777777778
$ws.Range("S15").Value = 0.2755555555555556
$ws.Range("F16").Value = 0.01565995525727069
$ws.Range("H16").Value = 0.2237136465324385
$ws.Range("I16").Value = 0.0738255033557047
$ws.Range("J16").Value = 0.3557046979865772
$ws.Range("K16").Value = 0.1319910514541387
$ws.Range("M16").Value = 0.01342281879194631
$ws.Range("N16").Value = 0.002237136465324385
$ws.Range("O16").Value = 0.03355704697986577
$ws.Range("S16").Value = 0.1498881431767338
$ws.Range("F17").Value = 0.01785714285714286
$ws.Range("H17").Value = 0.2119047619047619
$ws.Range("I17").Value = 0.07142857142857142
$ws.Range("J17").Value = 0.4178571428571429
$ws.Range("K17").Value = 0.1095238095238095
$ws.Range("M17").Value = 0.01547619047619048
$ws.Range("N17").Value = 0.001190476190476191
$ws.Range("O17").Value = 0.05238095238095238
$ws.Range("S17").Value = 0.1023809523809524
$ws.Range("F18").Value = 0.02127659574468085
$ws.Range("H18").Value = 0.2340425531914894
$ws.Range("I18").Value = 0.05319148936170213
$ws.Range("J18").Value = 0.4122340425531915
$ws.Range("K18").Value = 0.0797872340425532
$ws.Range("M18").Value = 0.01595744680851064
$ws.Range("N18").Value = 0.002659574468085106
$ws.Range("O18").Value = 0.05851063829787234
$ws.Range("S18").Value = 0.1223404255319149
$ws.Range("F19").Value = 0.0206794682422452
$ws.Range("H19").Value = 0.2208271787296898
$ws.Range("I19").Value = 0.07865583456425407
$ws.Range("J19").Value = 0.3607828655834564
$ws.Range("K19").Value = 0.1100443131462334
$ws.Range("M19").Value = 0.02363367799113737
$ws.Range("N19").Value = 0.001107828655834564
$ws.Range("O19").Value = 0.06794682422451995
$ws.Range("S19").Value = 0.1163220088626293
